{"js": "// Fix potential encoding problems in the document text:\n//  1. The heading \"Expression Replacement in global paragraphs\" was split\n//     across two runs; merge it back into a single run/text flow.\n//  2. Replace the German low-high quote marks (\u201e ... \u201c) that were wrapping\n//     the placeholder variable names \"name\" and \"foo\" with plain text\n//     (no quote marks), since they were a source of encoding trouble.\n\nconst body = context.document.body;\n\n// 1) Merge the heading's two runs into one by rewriting the whole\n//    paragraph range with its already-concatenated text.\nconst headingParagraph = body.paragraphs.getFirst();\nconst headingRange = headingParagraph.getRange();\nheadingRange.load(\"text\");\nawait context.sync();\nheadingRange.insertText(headingRange.text, \"Replace\");\nawait context.sync();\n\n// 2) Drop the \u201e \u201c quote marks around the variable names.\nconst nameHits = body.search(\"\\u201ename\\u201c\", { matchCase: true });\nnameHits.load(\"text\");\nawait context.sync();\nif (nameHits.items.length > 0) {\n  nameHits.items[0].insertText(\"name\", \"Replace\");\n}\n\nconst fooHits = body.search(\"\\u201efoo\\u201c\", { matchCase: true });\nfooHits.load(\"text\");\nawait context.sync();\nif (fooHits.items.length > 0) {\n  fooHits.items[0].insertText(\"foo\", \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Fix potential encoding problems in the document text:\n#  1. The heading \"Expression Replacement in global paragraphs\" was split\n#     across two runs; merge it back into a single run/text flow by running\n#     it through Find & Replace (re-writing the run it touches).\n#  2. Replace the German low-high quote marks (\u201e ... \u201c) that were wrapping\n#     the placeholder variable names \"name\" and \"foo\" with plain text\n#     (no quote marks), since they were a source of encoding trouble.\n\n$d = $word.ActiveDocument\n\n# 1) Merge the heading's two runs into one.\n$headingParagraph = $d.Paragraphs(1)\n$headingRange = $headingParagraph.Range\n$headingFind = $headingRange.Find\n$headingFind.ClearFormatting()\n$headingFind.Replacement.ClearFormatting()\n$headingFind.Execute(\"Expression\", $false, $false, $false, $false, $false, $true, 1, $false, \"Expression\", 2)\n\n# 2) Drop the \u201e \u201c quote marks around the variable names.\n$openQuote = [char]0x201E\n$closeQuote = [char]0x201C\n\n$nameRange = $d.Content\n$nameFind = $nameRange.Find\n$nameFind.ClearFormatting()\n$nameFind.Replacement.ClearFormatting()\n$nameFind.Execute($openQuote + \"name\" + $closeQuote, $false, $false, $false, $false, $false, $true, 1, $false, \"name\", 2)\n\n$fooRange = $d.Content\n$fooFind = $fooRange.Find\n$fooFind.ClearFormatting()\n$fooFind.Replacement.ClearFormatting()\n$fooFind.Execute($openQuote + \"foo\" + $closeQuote, $false, $false, $false, $false, $false, $true, 1, $false, \"foo\", 2)\n"}
